# Insert a new data row at row 166 (pushing existing rows 166:224 down to 167:225)
# and populate it with the new price-report entry for Berenjena.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("166:166").Insert()

$ws.Range("A166").Value = 3
$ws.Range("B166").Value = "Femacal de La Calera"
$ws.Range("C166").Value = "Coquimbo"
$ws.Range("D166").Value = 44588
$ws.Range("E166").Value = 5
$ws.Range("F166").Value = 100112001
$ws.Range("G166").Value = "Berenjena"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 115
$ws.Range("K166").Value = 9500
$ws.Range("L166").Value = 10000
$ws.Range("M166").Value = 9761
$ws.Range("N166").Value = "`$/caja 60 unidades"
$ws.Range("O166").Value = "Región de Arica y Parinacota"
$ws.Range("P166").Value = 163
$ws.Range("Q166").Value = 60
$ws.Range("R166").Value = "Hortaliza"
